$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1844933.5
$ws.Range("I9").Value = 645.75
$ws.Range("J9").Value = 2664617
$ws.Range("K9").Value = 645.75
$ws.Range("L9").Value = 2664617
$ws.Range("M9").Value = -476.75
$ws.Range("N9").Value = -2664955

$ws.Range("H29").Value = 939.3333
$ws.Range("J29").Value = 5003.5
$ws.Range("L29").Value = 15010.5
$ws.Range("N29").Value = -15572.5

$ws.Range("H33").Value = 297.69232
$ws.Range("I33").Value = 342.22223
$ws.Range("K33").Value = 342.22223
$ws.Range("M33").Value = -113.22223

$ws.Range("H38").Value = 127.85714
$ws.Range("I38").Value = 127.85714
$ws.Range("K38").Value = 383.57142
$ws.Range("M38").Value = -11.57141999999999

$ws.Range("H43").Value = 5838.826
$ws.Range("I43").Value = 2012
$ws.Range("J43").Value = 7879.8
$ws.Range("K43").Value = 2012
$ws.Range("L43").Value = 7879.8
$ws.Range("M43").Value = -1943
$ws.Range("N43").Value = -8017.8

$ws.Range("H48").Value = 1089
$ws.Range("I48").Value = 1089
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3267
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2975
$ws.Range("N48").ClearContents()

$ws.Range("H56").Value = 1089
$ws.Range("I56").Value = 1089
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3267
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2733
$ws.Range("N56").ClearContents()

$ws.Range("H80").Value = 1971.2
$ws.Range("I80").Value = 1096.5
$ws.Range("J80").Value = 2289.2727
$ws.Range("K80").Value = 3289.5
$ws.Range("L80").Value = 6867.8181
$ws.Range("M80").Value = -2291.5
$ws.Range("N80").Value = -8863.8181

$ws.Range("H83").Value = 1971.2
$ws.Range("I83").Value = 1096.5
$ws.Range("J83").Value = 2289.2727
$ws.Range("K83").Value = 9868.5
$ws.Range("L83").Value = 20603.4543
$ws.Range("M83").Value = -4876.5
$ws.Range("N83").Value = -30587.4543

$ws.Range("H86").Value = 2291.7778
$ws.Range("I86").Value = 2680
$ws.Range("J86").Value = 1806.5
$ws.Range("K86").Value = 2680
$ws.Range("L86").Value = 1806.5
$ws.Range("M86").Value = -1557
$ws.Range("N86").Value = -4052.5

$ws.Range("H89").Value = 2291.7778
$ws.Range("I89").Value = 2680
$ws.Range("J89").Value = 1806.5
$ws.Range("K89").Value = 13400
$ws.Range("L89").Value = 9032.5
$ws.Range("M89").Value = -7784
$ws.Range("N89").Value = -20264.5

$ws.Range("H92").Value = 1718.6072
$ws.Range("I92").Value = 1164.88
$ws.Range("J92").Value = 6333
$ws.Range("K92").Value = 1164.88
$ws.Range("L92").Value = 6333
$ws.Range("M92").Value = 83.11999999999989
$ws.Range("N92").Value = -8829

$ws.Range("H112").Value = 3248186.5
$ws.Range("J112").Value = 3789384.8
$ws.Range("L112").Value = 11368154.4
$ws.Range("N112").Value = -11370370.4

$ws.Range("H138").Value = 3690.762
$ws.Range("I138").Value = 3940.2222
$ws.Range("J138").Value = 3503.6667
$ws.Range("K138").Value = 11820.6666
$ws.Range("L138").Value = 10511.0001
$ws.Range("M138").Value = -6680.6666
$ws.Range("N138").Value = -20791.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33055.45
$ws.Range("I32").Value = 33364.562
$ws.Range("K32").Value = 33364.562
$ws.Range("M32").Value = -33077.562

$ws.Range("H45").Value = 1614.375
$ws.Range("I45").Value = 1527.7693
$ws.Range("J45").Value = 1989.6666
$ws.Range("K45").Value = 1527.7693
$ws.Range("L45").Value = 1989.6666
$ws.Range("M45").Value = -1150.7693
$ws.Range("N45").Value = -2743.6666

$ws.Range("H88").Value = 1331.6666
$ws.Range("I88").Value = 1585
$ws.Range("J88").Value = 1247.2222
$ws.Range("K88").Value = 1585
$ws.Range("L88").Value = 1247.2222
$ws.Range("M88").Value = -1179
$ws.Range("N88").Value = -2059.2222

$ws.Range("H91").Value = 1331.6666
$ws.Range("I91").Value = 1585
$ws.Range("J91").Value = 1247.2222
$ws.Range("K91").Value = 1585
$ws.Range("L91").Value = 1247.2222
$ws.Range("M91").Value = -181
$ws.Range("N91").Value = -4055.2222

$ws.Range("H102").Value = 16735167
$ws.Range("I102").Value = 2115.6924
$ws.Range("K102").Value = 2115.6924
$ws.Range("M102").Value = -493.6923999999999

$ws.Range("H110").Value = 17859412
$ws.Range("I110").Value = 62500324
$ws.Range("K110").Value = 62500324
$ws.Range("M110").Value = -62498279

$ws.Range("H139").Value = 99998
$ws.Range("J139").Value = 99998
$ws.Range("L139").Value = 99998
$ws.Range("N139").Value = -110278

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2320.7
$ws.Range("I94").Value = 2634.2222
$ws.Range("J94").Value = 2064.182
$ws.Range("K94").Value = 2634.2222
$ws.Range("L94").Value = 2064.182
$ws.Range("M94").Value = -2183.2222
$ws.Range("N94").Value = -2966.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2765.3333
$ws.Range("J16").Value = 2200
$ws.Range("L16").Value = 2200
$ws.Range("N16").Value = -2774

$ws.Range("H22").Value = 415.2
$ws.Range("I22").Value = 430.57144
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 430.57144
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -80.57144
$ws.Range("N22").Value = -900

$ws.Range("H31").Value = 5538.1304
$ws.Range("I31").Value = 2532.4167
$ws.Range("J31").Value = 8817.091
$ws.Range("K31").Value = 2532.4167
$ws.Range("L31").Value = 8817.091
$ws.Range("M31").Value = -2237.4167
$ws.Range("N31").Value = -9407.091

$ws.Range("H34").Value = 5538.1304
$ws.Range("I34").Value = 2532.4167
$ws.Range("J34").Value = 8817.091
$ws.Range("K34").Value = 2532.4167
$ws.Range("L34").Value = 8817.091
$ws.Range("M34").Value = -2330.4167
$ws.Range("N34").Value = -9221.091

$ws.Range("H93").Value = 20001
$ws.Range("I93").Value = 20001
$ws.Range("K93").Value = 20001
$ws.Range("M93").Value = -18129

$ws.Range("H113").Value = 2765.3333
$ws.Range("J113").Value = 2200
$ws.Range("L113").Value = 2200
$ws.Range("N113").Value = -6540

$ws.Range("H134").Value = 6406.55
$ws.Range("I134").Value = 5286.3125
$ws.Range("K134").Value = 15858.9375
$ws.Range("M134").Value = -13323.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 573.63635
$ws.Range("J92").Value = 858
$ws.Range("L92").Value = 2574
$ws.Range("N92").Value = -5070

$ws.Range("H113").Value = 2289.0454
$ws.Range("I113").Value = 3305.4
$ws.Range("J113").Value = 1990.1177
$ws.Range("K113").Value = 9916.200000000001
$ws.Range("L113").Value = 5970.3531
$ws.Range("M113").Value = -7746.200000000001
$ws.Range("N113").Value = -10310.3531

$ws.Range("H131").Value = 25647260
$ws.Range("I131").Value = 166667300
$ws.Range("J131").Value = 7252.909
$ws.Range("K131").Value = 500001900
$ws.Range("L131").Value = 21758.727
$ws.Range("M131").Value = -499996860
$ws.Range("N131").Value = -31838.727

$ws.Range("H132").Value = 54665.05
$ws.Range("J132").Value = 2430.4614
$ws.Range("L132").Value = 21874.1526
$ws.Range("N132").Value = -26934.1526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3516
$ws.Range("I70").Value = 4774
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 4774
$ws.Range("L70").Value = 1000
$ws.Range("M70").Value = -4504
$ws.Range("N70").Value = -1540

$ws.Range("H73").Value = 3516
$ws.Range("I73").Value = 4774
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 4774
$ws.Range("L73").Value = 1000
$ws.Range("M73").Value = -3838
$ws.Range("N73").Value = -2872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3792.3438
$ws.Range("I22").Value = 2665.0715
$ws.Range("K22").Value = 2665.0715
$ws.Range("M22").Value = -2370.0715

$ws.Range("H27").Value = 3792.3438
$ws.Range("I27").Value = 2665.0715
$ws.Range("K27").Value = 2665.0715
$ws.Range("M27").Value = -2558.0715

$ws.Range("H40").Value = 125003000
$ws.Range("I40").Value = 125003000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 125003000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -125002864
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 6285.2334
$ws.Range("I46").Value = 1233.6666
$ws.Range("K46").Value = 1233.6666
$ws.Range("M46").Value = -1045.6666

$ws.Range("H93").Value = 1253.1904
$ws.Range("I93").Value = 1151.2667
$ws.Range("J93").Value = 1508
$ws.Range("K93").Value = 1151.2667
$ws.Range("L93").Value = 1508
$ws.Range("M93").Value = 96.7333000000001
$ws.Range("N93").Value = -4004

$ws.Range("H132").Value = 5225
$ws.Range("I132").Value = 4269.16
$ws.Range("J132").Value = 6818.067
$ws.Range("K132").Value = 12807.48
$ws.Range("L132").Value = 20454.201
$ws.Range("M132").Value = -10277.48
$ws.Range("N132").Value = -25514.201

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 32399.2
$ws.Range("J101").Value = 32399.2
$ws.Range("L101").Value = 32399.2
$ws.Range("N101").Value = -38889.2

$ws.Range("H116").Value = 83149
$ws.Range("J116").Value = 83149
$ws.Range("L116").Value = 83149
$ws.Range("N116").Value = -92327

$ws.Range("H122").Value = 2306.9756
$ws.Range("I122").Value = 2083.7568
$ws.Range("J122").Value = 4371.75
$ws.Range("K122").Value = 6251.2704
$ws.Range("L122").Value = 13115.25
$ws.Range("M122").Value = -3801.2704
$ws.Range("N122").Value = -18015.25

$ws.Range("H132").Value = 6807.9688
$ws.Range("I132").Value = 5875.2085
$ws.Range("J132").Value = 9606.25
$ws.Range("K132").Value = 17625.6255
$ws.Range("L132").Value = 28818.75
$ws.Range("M132").Value = -15095.6255
$ws.Range("N132").Value = -33878.75
